$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 6 new rows (7-12) below row 6, copying row 6's border/format so the
#     existing middle-band style (thin border, no top) carries through ---
$ws.Range("A7:D12").EntireRow.Insert()
$ws.Range("A6:D6").Copy()
$ws.Range("A7:D12").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A7:D12").RowHeight = 30
$excel.CutCopyMode = 0

# --- Date of Defense value (B3), formatted as a short date (built-in date format) ---
$ws.Range("B3").Value2 = 42282
$ws.Range("B3").NumberFormat = "mm-dd-yy"

# --- Row 6: first panelist + first revision note, marked Done ---
$ws.Range("A6").Value = "Sir. Mike Dela Fuente"
$ws.Range("B6").Value = "Equation for Accuracy"
$ws.Range("D6").Value = "Done"

# --- Row 7 ---
$ws.Range("B7").Value = "Computation for Specificity"
$ws.Range("D7").Value = "Done"

# --- Row 8 ---
$ws.Range("B8").Value = "Discuss TP, FP, TN, TP"
$ws.Range("D8").Value = "Done"

# --- Row 9 ---
$ws.Range("B9").Value = "Use of Post-test and Pre-Test tanggalin na"
$ws.Range("D9").Value = "Done"

# --- Row 10 (still pending, no Status yet) ---
$ws.Range("B10").Value = "What will be your Population"

# --- Row 11 ---
$ws.Range("B11").Value = "Sampling Technique"

# --- Row 12 ---
$ws.Range("B12").Value = "Hypothesis ibase sa SOP"

# --- Row 13: second panelist + her revision note, marked Done ---
$ws.Range("A13").Value = "Maam Ria A. Sagum"
$ws.Range("B13").Value = "Tanggalin na number 2 sa SOP"
$ws.Range("D13").Value = "Done"

# --- Row 14 ---
$ws.Range("B14").Value = "Palitan yung number 1 sa SOP ng accuracy"

# --- Column B needs to be wide enough to fit the longest revision note ---
$ws.Range("B1").ColumnWidth = 37.91

# --- Leave the selection where the editor left off ---
$ws.Range("D10").Select()
